# Updates cryptos list: prices/volume percentages refreshed, and three
# coin pairs swapped position (rows 30-32, 39-40, 46-47) to reflect new ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) and E (Volume(1h)) columns hold number-like / percent-like text
# (e.g. "65.912.90", "  +2.48%  "). Force text entry via NumberFormat "@" so
# Excel does not coerce them to numeric values, then ClearFormats() so the
# cell keeps using the default style (matches the original, un-styled cells).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "65.912.90"
Set-TextValue "E2" "  +2.48%  "

# Row 3
Set-TextValue "D3" "3.219.30"
Set-TextValue "E3" "  +4.65%  "

# Row 4
Set-TextValue "E4" "  +0.21%  "

# Row 5
Set-TextValue "D5" "570.58"
Set-TextValue "E5" "  +1.04%  "

# Row 6
Set-TextValue "D6" "154.28"
Set-TextValue "E6" "  +8.46%  "

# Row 7
Set-TextValue "E7" "  +0.11%  "

# Row 8
Set-TextValue "D8" "3.207.20"
Set-TextValue "E8" "  +4.57%  "

# Row 9
Set-TextValue "D9" "0.510"
Set-TextValue "E9" "  +3.44%  "

# Row 10
Set-TextValue "D10" "7.17"
Set-TextValue "E10" "  +12.62%  "

# Row 11
Set-TextValue "D11" "0.165"
Set-TextValue "E11" "  +3.95%  "

# Row 12
Set-TextValue "D12" "0.480"
Set-TextValue "E12" "  +3.25%  "

# Row 13
Set-TextValue "D13" "37.63"
Set-TextValue "E13" "  +5.33%  "

# Row 14
Set-TextValue "D14" "0.0000232"
Set-TextValue "E14" "  +3.21%  "

# Row 15
Set-TextValue "D15" "3.714.52"
Set-TextValue "E15" "  +4.56%  "

# Row 16
Set-TextValue "D16" "65.957.60"
Set-TextValue "E16" "  +2.70%  "

# Row 17
Set-TextValue "D17" "541.75"
Set-TextValue "E17" "  +9.75%  "

# Row 18
Set-TextValue "E18" "  +2.65%  "

# Row 19
Set-TextValue "D19" "3.195.48"
Set-TextValue "E19" "  +4.24%  "

# Row 20
Set-TextValue "D20" "6.98"
Set-TextValue "E20" "  +4.77%  "

# Row 21
Set-TextValue "D21" "14.34"
Set-TextValue "E21" "  +4.29%  "

# Row 22
Set-TextValue "D22" "0.730"
Set-TextValue "E22" "  +6.02%  "

# Row 23
Set-TextValue "D23" "7.72"
Set-TextValue "E23" "  +6.68%  "

# Row 24
Set-TextValue "D24" "13.35"
Set-TextValue "E24" "  +5.85%  "

# Row 25
Set-TextValue "D25" "80.73"
Set-TextValue "E25" "  +3.13%  "

# Row 26
Set-TextValue "E26" "  +0.11%  "

# Row 27
Set-TextValue "D27" "9.34"
Set-TextValue "E27" "  +18.91%  "

# Row 28
Set-TextValue "D28" "2.86"
Set-TextValue "E28" "  +3.54%  "

# Row 29
Set-TextValue "D29" "2.24"
Set-TextValue "E29" "  +7.87%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D30" "27.48"
Set-TextValue "E30" "  +3.86%  "

# Row 31
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D31" "2.77"
Set-TextValue "E31" "  +5.25%  "

# Row 32
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D32" "0.998"
Set-TextValue "E32" "  -0.09%  "

# Row 33
Set-TextValue "D33" "1.18"
Set-TextValue "E33" "  +5.82%  "

# Row 34
Set-TextValue "D34" "568.58"
Set-TextValue "E34" "  +10.38%  "

# Row 35
Set-TextValue "D35" "5.73"
Set-TextValue "E35" "  +4.41%  "

# Row 36
Set-TextValue "D36" "6.31"
Set-TextValue "E36" "  +6.32%  "

# Row 37
Set-TextValue "D37" "0.0461"
Set-TextValue "E37" "  +13.94%  "

# Row 38
Set-TextValue "D38" "53.84"
Set-TextValue "E38" "  +1.22%  "

# Row 39
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D39" "3.07"
Set-TextValue "E39" "  +15.91%  "

# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D40" "0.0857"
Set-TextValue "E40" "  +7.44%  "

# Row 41
Set-TextValue "D41" "0.126"
Set-TextValue "E41" "  +4.48%  "

# Row 42
Set-TextValue "D42" "3.124.25"
Set-TextValue "E42" "  +6.56%  "

# Row 43
Set-TextValue "D43" "8.54"
Set-TextValue "E43" "  +2.41%  "

# Row 44
Set-TextValue "D44" "2.33"
Set-TextValue "E44" "  +11.37%  "

# Row 45
Set-TextValue "D45" "0.271"
Set-TextValue "E45" "  +10.45%  "

# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D46" "26.45"
Set-TextValue "E46" "  +5.94%  "

# Row 47
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D47" "1.00"
Set-TextValue "E47" "  +0.10%  "

# Row 48
Set-TextValue "D48" "0.0₃0553"
Set-TextValue "E48" "  +2.47%  "

# Row 49
Set-TextValue "D49" "0.112"
Set-TextValue "E49" "  +3.78%  "

# Row 50
Set-TextValue "D50" "121.84"
Set-TextValue "E50" "  +0.20%  "

# Row 51
Set-TextValue "D51" "2.21"
Set-TextValue "E51" "  +7.05%  "
